$d = $word.ActiveDocument

$d.Content.Find.Execute("3+37=40", $false, $false, $false, $false, $false, $true, 1, $false, "28+60=88", 2) | Out-Null
$d.Content.Find.Execute("69+24=93", $false, $false, $false, $false, $false, $true, 1, $false, "42+24=66", 2) | Out-Null
$d.Content.Find.Execute("81-75=6", $false, $false, $false, $false, $false, $true, 1, $false, "50+29=79", 2) | Out-Null
$d.Content.Find.Execute("71-6=65", $false, $false, $false, $false, $false, $true, 1, $false, "56+19=75", 2) | Out-Null
$d.Content.Find.Execute("52-10=42", $false, $false, $false, $false, $false, $true, 1, $false, "69+10=79", 2) | Out-Null
$d.Content.Find.Execute("94-21=73", $false, $false, $false, $false, $false, $true, 1, $false, "60-34=26", 2) | Out-Null
$d.Content.Find.Execute("77-58=19", $false, $false, $false, $false, $false, $true, 1, $false, "10+12=22", 2) | Out-Null
$d.Content.Find.Execute("49-44=5", $false, $false, $false, $false, $false, $true, 1, $false, "52+33=85", 2) | Out-Null
$d.Content.Find.Execute("59-35=24", $false, $false, $false, $false, $false, $true, 1, $false, "45+35=80", 2) | Out-Null
$d.Content.Find.Execute("49+2=51", $false, $false, $false, $false, $false, $true, 1, $false, "16-7=9", 2) | Out-Null
$d.Content.Find.Execute("41-9=32", $false, $false, $false, $false, $false, $true, 1, $false, "50+44=94", 2) | Out-Null
$d.Content.Find.Execute("43-14=29", $false, $false, $false, $false, $false, $true, 1, $false, "26+31=57", 2) | Out-Null
$d.Content.Find.Execute("33-28=5", $false, $false, $false, $false, $false, $true, 1, $false, "42+30=72", 2) | Out-Null
$d.Content.Find.Execute("36+53=89", $false, $false, $false, $false, $false, $true, 1, $false, "6+33=39", 2) | Out-Null
$d.Content.Find.Execute("75-70=5", $false, $false, $false, $false, $false, $true, 1, $false, "54-38=16", 2) | Out-Null
$d.Content.Find.Execute("33+7=40", $false, $false, $false, $false, $false, $true, 1, $false, "86+1=87", 2) | Out-Null
$d.Content.Find.Execute("10+68=78", $false, $false, $false, $false, $false, $true, 1, $false, "95-88=7", 2) | Out-Null
$d.Content.Find.Execute("17+33=50", $false, $false, $false, $false, $false, $true, 1, $false, "28+52=80", 2) | Out-Null
$d.Content.Find.Execute("88-68=20", $false, $false, $false, $false, $false, $true, 1, $false, "73-51=22", 2) | Out-Null
$d.Content.Find.Execute("24+12=36", $false, $false, $false, $false, $false, $true, 1, $false, "34+35=69", 2) | Out-Null
$d.Content.Find.Execute("41+29=70", $false, $false, $false, $false, $false, $true, 1, $false, "76+5=81", 2) | Out-Null
$d.Content.Find.Execute("63-59=4", $false, $false, $false, $false, $false, $true, 1, $false, "93+2=95", 2) | Out-Null
$d.Content.Find.Execute("78-57=21", $false, $false, $false, $false, $false, $true, 1, $false, "43+16=59", 2) | Out-Null
$d.Content.Find.Execute("86+12=98", $false, $false, $false, $false, $false, $true, 1, $false, "93-55=38", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $false, $false, $false, $false, $false, $true, 1, $false, "33+53=86", 2) | Out-Null
$d.Content.Find.Execute("55+6=61", $false, $false, $false, $false, $false, $true, 1, $false, "62-47=15", 2) | Out-Null
$d.Content.Find.Execute("11+67=78", $false, $false, $false, $false, $false, $true, 1, $false, "87-65=22", 2) | Out-Null
$d.Content.Find.Execute("48+12=60", $false, $false, $false, $false, $false, $true, 1, $false, "36+32=68", 2) | Out-Null
$d.Content.Find.Execute("10+78=88", $false, $false, $false, $false, $false, $true, 1, $false, "48-23=25", 2) | Out-Null
$d.Content.Find.Execute("95-65=30", $false, $false, $false, $false, $false, $true, 1, $false, "58-58=0", 2) | Out-Null
$d.Content.Find.Execute("73+12=85", $false, $false, $false, $false, $false, $true, 1, $false, "38+14=52", 2) | Out-Null
$d.Content.Find.Execute("40+21=61", $false, $false, $false, $false, $false, $true, 1, $false, "93-46=47", 2) | Out-Null
$d.Content.Find.Execute("82+13=95", $false, $false, $false, $false, $false, $true, 1, $false, "86-6=80", 2) | Out-Null
$d.Content.Find.Execute("36+0=36", $false, $false, $false, $false, $false, $true, 1, $false, "33+32=65", 2) | Out-Null
$d.Content.Find.Execute("28+18=46", $false, $false, $false, $false, $false, $true, 1, $false, "85+6=91", 2) | Out-Null
$d.Content.Find.Execute("97-46=51", $false, $false, $false, $false, $false, $true, 1, $false, "92-11=81", 2) | Out-Null
$d.Content.Find.Execute("1+4=5", $false, $false, $false, $false, $false, $true, 1, $false, "24-11=13", 2) | Out-Null
$d.Content.Find.Execute("46-45=1", $false, $false, $false, $false, $false, $true, 1, $false, "79-1=78", 2) | Out-Null
$d.Content.Find.Execute("0+13=13", $false, $false, $false, $false, $false, $true, 1, $false, "76-54=22", 2) | Out-Null
$d.Content.Find.Execute("86-85=1", $false, $false, $false, $false, $false, $true, 1, $false, "59-30=29", 2) | Out-Null
$d.Content.Find.Execute("89-1=88", $false, $false, $false, $false, $false, $true, 1, $false, "54+15=69", 2) | Out-Null
$d.Content.Find.Execute("43-22=21", $false, $false, $false, $false, $false, $true, 1, $false, "43+38=81", 2) | Out-Null
$d.Content.Find.Execute("25+64=89", $false, $false, $false, $false, $false, $true, 1, $false, "93-81=12", 2) | Out-Null
$d.Content.Find.Execute("99-44=55", $false, $false, $false, $false, $false, $true, 1, $false, "49-49=0", 2) | Out-Null
$d.Content.Find.Execute("48-30=18", $false, $false, $false, $false, $false, $true, 1, $false, "51+5=56", 2) | Out-Null
$d.Content.Find.Execute("50-31=19", $false, $false, $false, $false, $false, $true, 1, $false, "21+11=32", 2) | Out-Null
$d.Content.Find.Execute("2+79=81", $false, $false, $false, $false, $false, $true, 1, $false, "61-24=37", 2) | Out-Null
$d.Content.Find.Execute("67+3=70", $false, $false, $false, $false, $false, $true, 1, $false, "74+7=81", 2) | Out-Null
$d.Content.Find.Execute("98-94=4", $false, $false, $false, $false, $false, $true, 1, $false, "30+4=34", 2) | Out-Null
$d.Content.Find.Execute("5-1=4", $false, $false, $false, $false, $false, $true, 1, $false, "50-46=4", 2) | Out-Null
$d.Content.Find.Execute("92-40=52", $false, $false, $false, $false, $false, $true, 1, $false, "98-2=96", 2) | Out-Null
$d.Content.Find.Execute("29+55=84", $false, $false, $false, $false, $false, $true, 1, $false, "75-41=34", 2) | Out-Null
$d.Content.Find.Execute("87-24=63", $false, $false, $false, $false, $false, $true, 1, $false, "16+70=86", 2) | Out-Null
$d.Content.Find.Execute("32+53=85", $false, $false, $false, $false, $false, $true, 1, $false, "37-16=21", 2) | Out-Null
$d.Content.Find.Execute("22+59=81", $false, $false, $false, $false, $false, $true, 1, $false, "39+6=45", 2) | Out-Null
$d.Content.Find.Execute("76+23=99", $false, $false, $false, $false, $false, $true, 1, $false, "2+47=49", 2) | Out-Null
$d.Content.Find.Execute("12+58=70", $false, $false, $false, $false, $false, $true, 1, $false, "82-27=55", 2) | Out-Null
$d.Content.Find.Execute("5+23=28", $false, $false, $false, $false, $false, $true, 1, $false, "6+18=24", 2) | Out-Null
$d.Content.Find.Execute("87+0=87", $false, $false, $false, $false, $false, $true, 1, $false, "64-48=16", 2) | Out-Null
$d.Content.Find.Execute("60-58=2", $false, $false, $false, $false, $false, $true, 1, $false, "15+45=60", 2) | Out-Null
$d.Content.Find.Execute("5+11=16", $false, $false, $false, $false, $false, $true, 1, $false, "55-33=22", 2) | Out-Null
$d.Content.Find.Execute("36-7=29", $false, $false, $false, $false, $false, $true, 1, $false, "18+43=61", 2) | Out-Null
$d.Content.Find.Execute("14+42=56", $false, $false, $false, $false, $false, $true, 1, $false, "98-24=74", 2) | Out-Null
$d.Content.Find.Execute("57+32=89", $false, $false, $false, $false, $false, $true, 1, $false, "52+20=72", 2) | Out-Null
$d.Content.Find.Execute("97-15=82", $false, $false, $false, $false, $false, $true, 1, $false, "85-7=78", 2) | Out-Null
$d.Content.Find.Execute("32+0=32", $false, $false, $false, $false, $false, $true, 1, $false, "89+2=91", 2) | Out-Null
$d.Content.Find.Execute("63+34=97", $false, $false, $false, $false, $false, $true, 1, $false, "98+1=99", 2) | Out-Null
$d.Content.Find.Execute("62+10=72", $false, $false, $false, $false, $false, $true, 1, $false, "1+54=55", 2) | Out-Null
$d.Content.Find.Execute("58-8=50", $false, $false, $false, $false, $false, $true, 1, $false, "54-13=41", 2) | Out-Null
$d.Content.Find.Execute("96-68=28", $false, $false, $false, $false, $false, $true, 1, $false, "37+47=84", 2) | Out-Null
$d.Content.Find.Execute("75-25=50", $false, $false, $false, $false, $false, $true, 1, $false, "54-42=12", 2) | Out-Null
$d.Content.Find.Execute("19+27=46", $false, $false, $false, $false, $false, $true, 1, $false, "29-22=7", 2) | Out-Null
$d.Content.Find.Execute("81-17=64", $false, $false, $false, $false, $false, $true, 1, $false, "35-6=29", 2) | Out-Null
$d.Content.Find.Execute("39-8=31", $false, $false, $false, $false, $false, $true, 1, $false, "34-32=2", 2) | Out-Null
$d.Content.Find.Execute("81+1=82", $false, $false, $false, $false, $false, $true, 1, $false, "24+55=79", 2) | Out-Null
$d.Content.Find.Execute("98-4=94", $false, $false, $false, $false, $false, $true, 1, $false, "42+35=77", 2) | Out-Null
$d.Content.Find.Execute("93-92=1", $false, $false, $false, $false, $false, $true, 1, $false, "24+27=51", 2) | Out-Null
$d.Content.Find.Execute("24+66=90", $false, $false, $false, $false, $false, $true, 1, $false, "49-12=37", 2) | Out-Null
$d.Content.Find.Execute("86-0=86", $false, $false, $false, $false, $false, $true, 1, $false, "85-58=27", 2) | Out-Null
$d.Content.Find.Execute("29+15=44", $false, $false, $false, $false, $false, $true, 1, $false, "44+13=57", 2) | Out-Null
$d.Content.Find.Execute("17+39=56", $false, $false, $false, $false, $false, $true, 1, $false, "42-21=21", 2) | Out-Null
$d.Content.Find.Execute("78-5=73", $false, $false, $false, $false, $false, $true, 1, $false, "21+78=99", 2) | Out-Null
$d.Content.Find.Execute("27+59=86", $false, $false, $false, $false, $false, $true, 1, $false, "86-8=78", 2) | Out-Null
$d.Content.Find.Execute("53-1=52", $false, $false, $false, $false, $false, $true, 1, $false, "27-15=12", 2) | Out-Null
$d.Content.Find.Execute("42-39=3", $false, $false, $false, $false, $false, $true, 1, $false, "56+2=58", 2) | Out-Null
$d.Content.Find.Execute("52+40=92", $false, $false, $false, $false, $false, $true, 1, $false, "98-36=62", 2) | Out-Null
$d.Content.Find.Execute("25+65=90", $false, $false, $false, $false, $false, $true, 1, $false, "44-36=8", 2) | Out-Null
$d.Content.Find.Execute("43-2=41", $false, $false, $false, $false, $false, $true, 1, $false, "61+4=65", 2) | Out-Null
$d.Content.Find.Execute("71-60=11", $false, $false, $false, $false, $false, $true, 1, $false, "33+59=92", 2) | Out-Null
$d.Content.Find.Execute("41+9=50", $false, $false, $false, $false, $false, $true, 1, $false, "83+14=97", 2) | Out-Null
$d.Content.Find.Execute("92-23=69", $false, $false, $false, $false, $false, $true, 1, $false, "66-7=59", 2) | Out-Null
$d.Content.Find.Execute("46+12=58", $false, $false, $false, $false, $false, $true, 1, $false, "76-32=44", 2) | Out-Null
$d.Content.Find.Execute("11+76=87", $false, $false, $false, $false, $false, $true, 1, $false, "45+9=54", 2) | Out-Null
$d.Content.Find.Execute("61-16=45", $false, $false, $false, $false, $false, $true, 1, $false, "9+2=11", 2) | Out-Null
$d.Content.Find.Execute("70-53=17", $false, $false, $false, $false, $false, $true, 1, $false, "34+39=73", 2) | Out-Null
$d.Content.Find.Execute("25+21=46", $false, $false, $false, $false, $false, $true, 1, $false, "45+45=90", 2) | Out-Null
$d.Content.Find.Execute("70+22=92", $false, $false, $false, $false, $false, $true, 1, $false, "30+60=90", 2) | Out-Null
$d.Content.Find.Execute("3+32=35", $false, $false, $false, $false, $false, $true, 1, $false, "83-48=35", 2) | Out-Null
$d.Content.Find.Execute("43-8=35", $false, $false, $false, $false, $false, $true, 1, $false, "93-18=75", 2) | Out-Null
$d.Content.Find.Execute("86-78=8", $false, $false, $false, $false, $false, $true, 1, $false, "0+88=88", 2) | Out-Null
